# Applies the diff: values in columns D, L, M, N, O, P, R, S were
# redistributed across rows 2,5-42,44-46 (a row-wise permutation of the
# underlying daily price records). This script writes the final target
# value for every affected cell directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = 44181
$ws.Cells.Item(2, 12).Value2 = 'Primera'
$ws.Cells.Item(2, 13).Value2 = 140
$ws.Cells.Item(2, 14).Value2 = 4000
$ws.Cells.Item(2, 15).Value2 = 4500
$ws.Cells.Item(2, 16).Value2 = 4250
$ws.Cells.Item(2, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(2, 19).Value2 = 2125

$ws.Cells.Item(5, 4).Value2 = 44194
$ws.Cells.Item(5, 12).Value2 = 'Primera'
$ws.Cells.Item(5, 13).Value2 = 250
$ws.Cells.Item(5, 14).Value2 = 4000
$ws.Cells.Item(5, 15).Value2 = 4000
$ws.Cells.Item(5, 16).Value2 = 4000
$ws.Cells.Item(5, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(5, 19).Value2 = 2000

$ws.Cells.Item(6, 4).Value2 = 44186
$ws.Cells.Item(6, 12).Value2 = 'Primera'
$ws.Cells.Item(6, 13).Value2 = 200
$ws.Cells.Item(6, 14).Value2 = 4000
$ws.Cells.Item(6, 15).Value2 = 4000
$ws.Cells.Item(6, 16).Value2 = 4000
$ws.Cells.Item(6, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(6, 19).Value2 = 2000

$ws.Cells.Item(7, 4).Value2 = 44189
$ws.Cells.Item(7, 12).Value2 = 'Primera'
$ws.Cells.Item(7, 13).Value2 = 300
$ws.Cells.Item(7, 14).Value2 = 3000
$ws.Cells.Item(7, 15).Value2 = 3000
$ws.Cells.Item(7, 16).Value2 = 3000
$ws.Cells.Item(7, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(7, 19).Value2 = 1500

$ws.Cells.Item(8, 4).Value2 = 44189
$ws.Cells.Item(8, 12).Value2 = 'Primera'
$ws.Cells.Item(8, 13).Value2 = 250
$ws.Cells.Item(8, 14).Value2 = 3000
$ws.Cells.Item(8, 15).Value2 = 3000
$ws.Cells.Item(8, 16).Value2 = 3000
$ws.Cells.Item(8, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(8, 19).Value2 = 1500

$ws.Cells.Item(9, 4).Value2 = 44202
$ws.Cells.Item(9, 12).Value2 = 'Primera'
$ws.Cells.Item(9, 13).Value2 = 200
$ws.Cells.Item(9, 14).Value2 = 4000
$ws.Cells.Item(9, 15).Value2 = 4000
$ws.Cells.Item(9, 16).Value2 = 4000
$ws.Cells.Item(9, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(9, 19).Value2 = 2000

$ws.Cells.Item(10, 4).Value2 = 44222
$ws.Cells.Item(10, 12).Value2 = 'Primera'
$ws.Cells.Item(10, 13).Value2 = 250
$ws.Cells.Item(10, 14).Value2 = 4000
$ws.Cells.Item(10, 15).Value2 = 4000
$ws.Cells.Item(10, 16).Value2 = 4000
$ws.Cells.Item(10, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(10, 19).Value2 = 2000

$ws.Cells.Item(11, 4).Value2 = 44222
$ws.Cells.Item(11, 12).Value2 = 'Primera'
$ws.Cells.Item(11, 13).Value2 = 300
$ws.Cells.Item(11, 14).Value2 = 4000
$ws.Cells.Item(11, 15).Value2 = 4000
$ws.Cells.Item(11, 16).Value2 = 4000
$ws.Cells.Item(11, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(11, 19).Value2 = 2000

$ws.Cells.Item(12, 4).Value2 = 44224
$ws.Cells.Item(12, 12).Value2 = 'Primera'
$ws.Cells.Item(12, 13).Value2 = 250
$ws.Cells.Item(12, 14).Value2 = 4000
$ws.Cells.Item(12, 15).Value2 = 4000
$ws.Cells.Item(12, 16).Value2 = 4000
$ws.Cells.Item(12, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(12, 19).Value2 = 2000

$ws.Cells.Item(13, 4).Value2 = 44224
$ws.Cells.Item(13, 12).Value2 = 'Primera'
$ws.Cells.Item(13, 13).Value2 = 300
$ws.Cells.Item(13, 14).Value2 = 4000
$ws.Cells.Item(13, 15).Value2 = 4000
$ws.Cells.Item(13, 16).Value2 = 4000
$ws.Cells.Item(13, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(13, 19).Value2 = 2000

$ws.Cells.Item(14, 4).Value2 = 44195
$ws.Cells.Item(14, 12).Value2 = 'Primera'
$ws.Cells.Item(14, 13).Value2 = 300
$ws.Cells.Item(14, 14).Value2 = 3000
$ws.Cells.Item(14, 15).Value2 = 3000
$ws.Cells.Item(14, 16).Value2 = 3000
$ws.Cells.Item(14, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(14, 19).Value2 = 1500

$ws.Cells.Item(15, 4).Value2 = 44169
$ws.Cells.Item(15, 12).Value2 = 'Primera'
$ws.Cells.Item(15, 13).Value2 = 200
$ws.Cells.Item(15, 14).Value2 = 5000
$ws.Cells.Item(15, 15).Value2 = 5000
$ws.Cells.Item(15, 16).Value2 = 5000
$ws.Cells.Item(15, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(15, 19).Value2 = 2500

$ws.Cells.Item(16, 4).Value2 = 44188
$ws.Cells.Item(16, 12).Value2 = 'Primera'
$ws.Cells.Item(16, 13).Value2 = 300
$ws.Cells.Item(16, 14).Value2 = 4000
$ws.Cells.Item(16, 15).Value2 = 4000
$ws.Cells.Item(16, 16).Value2 = 4000
$ws.Cells.Item(16, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(16, 19).Value2 = 2000

$ws.Cells.Item(17, 4).Value2 = 44188
$ws.Cells.Item(17, 12).Value2 = 'Primera'
$ws.Cells.Item(17, 13).Value2 = 500
$ws.Cells.Item(17, 14).Value2 = 4000
$ws.Cells.Item(17, 15).Value2 = 4000
$ws.Cells.Item(17, 16).Value2 = 4000
$ws.Cells.Item(17, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(17, 19).Value2 = 2000

$ws.Cells.Item(18, 4).Value2 = 44216
$ws.Cells.Item(18, 12).Value2 = 'Primera'
$ws.Cells.Item(18, 13).Value2 = 200
$ws.Cells.Item(18, 14).Value2 = 4000
$ws.Cells.Item(18, 15).Value2 = 4000
$ws.Cells.Item(18, 16).Value2 = 4000
$ws.Cells.Item(18, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(18, 19).Value2 = 2000

$ws.Cells.Item(19, 4).Value2 = 44216
$ws.Cells.Item(19, 12).Value2 = 'Primera'
$ws.Cells.Item(19, 13).Value2 = 400
$ws.Cells.Item(19, 14).Value2 = 4000
$ws.Cells.Item(19, 15).Value2 = 4000
$ws.Cells.Item(19, 16).Value2 = 4000
$ws.Cells.Item(19, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(19, 19).Value2 = 2000

$ws.Cells.Item(20, 4).Value2 = 44253
$ws.Cells.Item(20, 12).Value2 = 'Primera'
$ws.Cells.Item(20, 13).Value2 = 25
$ws.Cells.Item(20, 14).Value2 = 4000
$ws.Cells.Item(20, 15).Value2 = 4000
$ws.Cells.Item(20, 16).Value2 = 4000
$ws.Cells.Item(20, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(20, 19).Value2 = 2000

$ws.Cells.Item(21, 4).Value2 = 44250
$ws.Cells.Item(21, 12).Value2 = 'Primera'
$ws.Cells.Item(21, 13).Value2 = 100
$ws.Cells.Item(21, 14).Value2 = 4000
$ws.Cells.Item(21, 15).Value2 = 4000
$ws.Cells.Item(21, 16).Value2 = 4000
$ws.Cells.Item(21, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(21, 19).Value2 = 2000

$ws.Cells.Item(22, 4).Value2 = 44225
$ws.Cells.Item(22, 12).Value2 = 'Primera'
$ws.Cells.Item(22, 13).Value2 = 150
$ws.Cells.Item(22, 14).Value2 = 4000
$ws.Cells.Item(22, 15).Value2 = 4000
$ws.Cells.Item(22, 16).Value2 = 4000
$ws.Cells.Item(22, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(22, 19).Value2 = 2000

$ws.Cells.Item(23, 4).Value2 = 44225
$ws.Cells.Item(23, 12).Value2 = 'Primera'
$ws.Cells.Item(23, 13).Value2 = 200
$ws.Cells.Item(23, 14).Value2 = 4000
$ws.Cells.Item(23, 15).Value2 = 4000
$ws.Cells.Item(23, 16).Value2 = 4000
$ws.Cells.Item(23, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(23, 19).Value2 = 2000

$ws.Cells.Item(24, 4).Value2 = 44223
$ws.Cells.Item(24, 12).Value2 = 'Primera'
$ws.Cells.Item(24, 13).Value2 = 200
$ws.Cells.Item(24, 14).Value2 = 4000
$ws.Cells.Item(24, 15).Value2 = 4000
$ws.Cells.Item(24, 16).Value2 = 4000
$ws.Cells.Item(24, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(24, 19).Value2 = 2000

$ws.Cells.Item(25, 4).Value2 = 44193
$ws.Cells.Item(25, 12).Value2 = 'Primera'
$ws.Cells.Item(25, 13).Value2 = 200
$ws.Cells.Item(25, 14).Value2 = 3000
$ws.Cells.Item(25, 15).Value2 = 3000
$ws.Cells.Item(25, 16).Value2 = 3000
$ws.Cells.Item(25, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(25, 19).Value2 = 1500

$ws.Cells.Item(26, 4).Value2 = 44210
$ws.Cells.Item(26, 12).Value2 = 'Primera'
$ws.Cells.Item(26, 13).Value2 = 400
$ws.Cells.Item(26, 14).Value2 = 3000
$ws.Cells.Item(26, 15).Value2 = 4000
$ws.Cells.Item(26, 16).Value2 = 3500
$ws.Cells.Item(26, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(26, 19).Value2 = 1750

$ws.Cells.Item(27, 4).Value2 = 44217
$ws.Cells.Item(27, 12).Value2 = 'Primera'
$ws.Cells.Item(27, 13).Value2 = 250
$ws.Cells.Item(27, 14).Value2 = 4000
$ws.Cells.Item(27, 15).Value2 = 4000
$ws.Cells.Item(27, 16).Value2 = 4000
$ws.Cells.Item(27, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(27, 19).Value2 = 2000

$ws.Cells.Item(28, 4).Value2 = 44217
$ws.Cells.Item(28, 12).Value2 = 'Primera'
$ws.Cells.Item(28, 13).Value2 = 300
$ws.Cells.Item(28, 14).Value2 = 4000
$ws.Cells.Item(28, 15).Value2 = 4000
$ws.Cells.Item(28, 16).Value2 = 4000
$ws.Cells.Item(28, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(28, 19).Value2 = 2000

$ws.Cells.Item(29, 4).Value2 = 44201
$ws.Cells.Item(29, 12).Value2 = 'Primera'
$ws.Cells.Item(29, 13).Value2 = 200
$ws.Cells.Item(29, 14).Value2 = 4000
$ws.Cells.Item(29, 15).Value2 = 4000
$ws.Cells.Item(29, 16).Value2 = 4000
$ws.Cells.Item(29, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(29, 19).Value2 = 2000

$ws.Cells.Item(30, 4).Value2 = 44260
$ws.Cells.Item(30, 12).Value2 = 'Primera'
$ws.Cells.Item(30, 13).Value2 = 75
$ws.Cells.Item(30, 14).Value2 = 4000
$ws.Cells.Item(30, 15).Value2 = 4000
$ws.Cells.Item(30, 16).Value2 = 4000
$ws.Cells.Item(30, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(30, 19).Value2 = 2000

$ws.Cells.Item(31, 4).Value2 = 44187
$ws.Cells.Item(31, 12).Value2 = 'Primera'
$ws.Cells.Item(31, 13).Value2 = 100
$ws.Cells.Item(31, 14).Value2 = 3400
$ws.Cells.Item(31, 15).Value2 = 3400
$ws.Cells.Item(31, 16).Value2 = 3400
$ws.Cells.Item(31, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(31, 19).Value2 = 1700

$ws.Cells.Item(32, 4).Value2 = 44187
$ws.Cells.Item(32, 12).Value2 = 'Primera'
$ws.Cells.Item(32, 13).Value2 = 200
$ws.Cells.Item(32, 14).Value2 = 4000
$ws.Cells.Item(32, 15).Value2 = 4000
$ws.Cells.Item(32, 16).Value2 = 4000
$ws.Cells.Item(32, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(32, 19).Value2 = 2000

$ws.Cells.Item(33, 4).Value2 = 44187
$ws.Cells.Item(33, 12).Value2 = 'Segunda'
$ws.Cells.Item(33, 13).Value2 = 50
$ws.Cells.Item(33, 14).Value2 = 3000
$ws.Cells.Item(33, 15).Value2 = 3000
$ws.Cells.Item(33, 16).Value2 = 3000
$ws.Cells.Item(33, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(33, 19).Value2 = 1500

$ws.Cells.Item(34, 4).Value2 = 44196
$ws.Cells.Item(34, 12).Value2 = 'Primera'
$ws.Cells.Item(34, 13).Value2 = 150
$ws.Cells.Item(34, 14).Value2 = 4000
$ws.Cells.Item(34, 15).Value2 = 4000
$ws.Cells.Item(34, 16).Value2 = 4000
$ws.Cells.Item(34, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(34, 19).Value2 = 2000

$ws.Cells.Item(35, 4).Value2 = 44251
$ws.Cells.Item(35, 12).Value2 = 'Primera'
$ws.Cells.Item(35, 13).Value2 = 125
$ws.Cells.Item(35, 14).Value2 = 4000
$ws.Cells.Item(35, 15).Value2 = 4000
$ws.Cells.Item(35, 16).Value2 = 4000
$ws.Cells.Item(35, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(35, 19).Value2 = 2000

$ws.Cells.Item(36, 4).Value2 = 44252
$ws.Cells.Item(36, 12).Value2 = 'Primera'
$ws.Cells.Item(36, 13).Value2 = 75
$ws.Cells.Item(36, 14).Value2 = 4000
$ws.Cells.Item(36, 15).Value2 = 4000
$ws.Cells.Item(36, 16).Value2 = 4000
$ws.Cells.Item(36, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(36, 19).Value2 = 2000

$ws.Cells.Item(37, 4).Value2 = 44221
$ws.Cells.Item(37, 12).Value2 = 'Primera'
$ws.Cells.Item(37, 13).Value2 = 150
$ws.Cells.Item(37, 14).Value2 = 4000
$ws.Cells.Item(37, 15).Value2 = 4000
$ws.Cells.Item(37, 16).Value2 = 4000
$ws.Cells.Item(37, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(37, 19).Value2 = 2000

$ws.Cells.Item(38, 4).Value2 = 44221
$ws.Cells.Item(38, 12).Value2 = 'Primera'
$ws.Cells.Item(38, 13).Value2 = 200
$ws.Cells.Item(38, 14).Value2 = 4000
$ws.Cells.Item(38, 15).Value2 = 4000
$ws.Cells.Item(38, 16).Value2 = 4000
$ws.Cells.Item(38, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(38, 19).Value2 = 2000

$ws.Cells.Item(39, 4).Value2 = 44209
$ws.Cells.Item(39, 12).Value2 = 'Primera'
$ws.Cells.Item(39, 13).Value2 = 170
$ws.Cells.Item(39, 14).Value2 = 3000
$ws.Cells.Item(39, 15).Value2 = 4000
$ws.Cells.Item(39, 16).Value2 = 3500
$ws.Cells.Item(39, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(39, 19).Value2 = 1750

$ws.Cells.Item(40, 4).Value2 = 44215
$ws.Cells.Item(40, 12).Value2 = 'Primera'
$ws.Cells.Item(40, 13).Value2 = 750
$ws.Cells.Item(40, 14).Value2 = 4000
$ws.Cells.Item(40, 15).Value2 = 4000
$ws.Cells.Item(40, 16).Value2 = 4000
$ws.Cells.Item(40, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(40, 19).Value2 = 2000

$ws.Cells.Item(41, 4).Value2 = 44175
$ws.Cells.Item(41, 12).Value2 = 'Primera'
$ws.Cells.Item(41, 13).Value2 = 250
$ws.Cells.Item(41, 14).Value2 = 4000
$ws.Cells.Item(41, 15).Value2 = 4000
$ws.Cells.Item(41, 16).Value2 = 4000
$ws.Cells.Item(41, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(41, 19).Value2 = 2000

$ws.Cells.Item(42, 4).Value2 = 44203
$ws.Cells.Item(42, 12).Value2 = 'Primera'
$ws.Cells.Item(42, 13).Value2 = 350
$ws.Cells.Item(42, 14).Value2 = 4000
$ws.Cells.Item(42, 15).Value2 = 4000
$ws.Cells.Item(42, 16).Value2 = 4000
$ws.Cells.Item(42, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(42, 19).Value2 = 2000

$ws.Cells.Item(44, 4).Value2 = 44176
$ws.Cells.Item(44, 12).Value2 = 'Primera'
$ws.Cells.Item(44, 13).Value2 = 100
$ws.Cells.Item(44, 14).Value2 = 4000
$ws.Cells.Item(44, 15).Value2 = 4000
$ws.Cells.Item(44, 16).Value2 = 4000
$ws.Cells.Item(44, 18).Value2 = 'Provincia de Linares'
$ws.Cells.Item(44, 19).Value2 = 2000

$ws.Cells.Item(45, 4).Value2 = 44239
$ws.Cells.Item(45, 12).Value2 = 'Primera'
$ws.Cells.Item(45, 13).Value2 = 350
$ws.Cells.Item(45, 14).Value2 = 3500
$ws.Cells.Item(45, 15).Value2 = 4000
$ws.Cells.Item(45, 16).Value2 = 3750
$ws.Cells.Item(45, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(45, 19).Value2 = 1875

$ws.Cells.Item(46, 4).Value2 = 44211
$ws.Cells.Item(46, 12).Value2 = 'Primera'
$ws.Cells.Item(46, 13).Value2 = 200
$ws.Cells.Item(46, 14).Value2 = 3000
$ws.Cells.Item(46, 15).Value2 = 3500
$ws.Cells.Item(46, 16).Value2 = 3250
$ws.Cells.Item(46, 18).Value2 = 'Provincia de Curicó'
$ws.Cells.Item(46, 19).Value2 = 1625
